$d = $word.ActiveDocument

# Find the "Requisitos" bullet-list paragraph that contains the LOM3246 line
# (order there currently is: LOB1021, LOM3016, LOM3246). The LOM3246 line
# needs to move to the front of that list, ahead of LOB1021, while leaving
# the other two lines (and their line breaks) untouched.
$targetParagraph = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*LOM3246 -  Técnicas de Caracterização de Materiais  (Indicação de Conjunto)*") {
        $targetParagraph = $p
    }
}

if ($targetParagraph -ne $null) {
    $paragraphRange = $targetParagraph.Range

    # Locate the LOM3246 line within the paragraph.
    $lineRange = $paragraphRange.Duplicate
    $found = $lineRange.Find.Execute("LOM3246 -  Técnicas de Caracterização de Materiais  (Indicação de Conjunto)")

    if ($found) {
        # Extend the match to also capture the manual line break (w:br) that
        # immediately follows the line's text, so the whole "run" moves as a unit.
        [void]$lineRange.MoveEnd(1, 1)
        $movedText = $lineRange.Text

        # Remove the line from its current (last) position in the list.
        $lineRange.Delete()

        # Re-insert the captured text (with its trailing break) at the very
        # start of the paragraph, ahead of the "LOB1021" line.
        $insertionPoint = $targetParagraph.Range
        $insertionPoint.Collapse(1)
        $insertionPoint.InsertBefore($movedText)
    }
}
